$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '37.711.54'
$ws.Range('E2').Value = '  -0.12%  '
# Row 3
$ws.Range('D3').Value = '2.072.73'
$ws.Range('E3').Value = '  -1.68%  '
# Row 4
$ws.Range('E4').Value = '  +0.02%  '
# Row 5
$ws.Range('E5').Value = '  -0.48%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.622'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.00%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.39'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.16%  '
# Row 8
$ws.Range('E8').Value = '  +0.00%  '
# Row 9
$ws.Range('E9').Value = '  +1.11%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0782'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.69%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.107'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.34%  '
# Row 12
$ws.Range('D12').Value = '2.379.51'
$ws.Range('E12').Value = '  -1.58%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.72'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.94%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.91'
$ws.Range('D14').Style = 'Normal'
# Row 15
$ws.Range('E15').Value = '  -0.93%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.35'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.94%  '
# Row 17
$ws.Range('D17').Value = '2.088.57'
$ws.Range('E17').Value = '  -2.46%  '
# Row 18
$ws.Range('D18').Value = '37.646.23'
$ws.Range('E18').Value = '  -0.09%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.29%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.05%  '
# Row 21
$ws.Range('E21').Value = '  +1.61%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.99'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.41%  '
# Row 23
$ws.Range('E23').Value = '  -0.07%  '
# Row 24
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.30%  '
# Row 25
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.36'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.09%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.05'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.95%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.137'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.81%  '
# Row 28
$ws.Range('E28').Value = '  +0.49%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.44'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.09%  '
# Row 30
$ws.Range('E30').Value = '  -2.27%  '
# Row 31
$ws.Range('E31').Value = '  +2.62%  '
# Row 32
$ws.Range('E32').Value = '  +1.30%  '
# Row 33
$ws.Range('E33').Value = '  +1.64%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.64'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.17%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.48'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.66%  '
# Row 36
$ws.Range('E36').Value = '  +0.39%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.38'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.48%  '
# Row 38
$ws.Range('E38').Value = '  +0.02%  '
# Row 39
$ws.Range('E39').Value = '  -2.18%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '100.47'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.07%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0972'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.31%  '
# Row 43
$ws.Range('E43').Value = '  +0.78%  '
# Row 44
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.57'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.71%  '
# Row 45
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.442.79'
$ws.Range('E45').Value = '  -1.08%  '
# Row 46
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.15'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.96%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.19'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.66%  '
# Row 48
$ws.Range('E48').Value = '  +0.28%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.41'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.77%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.99'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.63%  '
# Row 51
$ws.Range('D51').Value = '2.264.92'
$ws.Range('E51').Value = '  -1.65%  '
